# Insert a new weekly price record at row 288 for
# "Feria Lagunitas de Puerto Montt - Brócoli", pushing the existing
# records (previously rows 288-387) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 288 (shifts 288..387 -> 289..388)
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new record's data
$ws.Range("A288").Value = 4
$ws.Range("B288").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C288").Value = "Los Lagos"
$ws.Range("D288").Value = "8/22/2022"
$ws.Range("E288").Value = 10
$ws.Range("F288").Value = 100112023
$ws.Range("G288").Value = "Brócoli"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 750
$ws.Range("K288").Value = 1500
$ws.Range("L288").Value = 1500
$ws.Range("M288").Value = 1500
$ws.Range("N288").Value = "`$/unidad"
$ws.Range("O288").Value = "Región del Maule"
$ws.Range("P288").Value = 1500
$ws.Range("Q288").Value = 1
$ws.Range("R288").Value = "Hortaliza"
